$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Prime the three brand-new cell formats in the same order Excel originally
# allocated them (cellXfs indices 8, 9, 10), so the xf table lines up byte
# for byte with the authored workbook.
# ---------------------------------------------------------------------------

# index 8: bordered, vertical-top only (used by the empty Metrics cells E12/E13)
$ws.Range("A9").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").HorizontalAlignment = -4142
$ws.Range("E12").VerticalAlignment = -4160
$ws.Range("E12").WrapText = $false

# index 9: bordered, left + vertical-top + wrap (used by Title/Authors/Link on row 11)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").HorizontalAlignment = -4131
$ws.Range("A11").VerticalAlignment = -4160
$ws.Range("A11").WrapText = $true

# index 10: bordered, left + vertical-top, no wrap (used by the empty Metrics cell E11)
$ws.Range("A9").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").HorizontalAlignment = -4131
$ws.Range("E11").VerticalAlignment = -4160
$ws.Range("E11").WrapText = $false

# ---------------------------------------------------------------------------
# Row 11 (ht=60)
# ---------------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").HorizontalAlignment = -4131
$ws.Range("B11").VerticalAlignment = -4160
$ws.Range("B11").WrapText = $true

$ws.Range("C9").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("A9").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").HorizontalAlignment = -4131
$ws.Range("D11").VerticalAlignment = -4160
$ws.Range("D11").WrapText = $true

# Values, set in the same order they first entered the shared-string table
$ws.Range("D11").Value = "https://www.dmst.aueb.gr/dds/pubs/conf/2008-OSS-qmodel/html/SGSS08.htm"
$ws.Range("A11").Value = "The SQO-OSS quality model: measurement based open source software evaluation"
$ws.Range("B11").Value = "Ioannis Samoladas, Georgios Gousios, Diomidis Spinellis and Ioannis Stamelos"
$ws.Range("C11").Value = 39692

$ws.Rows.Item(11).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 12 (ht=45)
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("A5").Copy()
$ws.Range("B12").PasteSpecial(-4122)

$ws.Range("C9").Copy()
$ws.Range("C12").PasteSpecial(-4122)

$ws.Range("A7").Copy()
$ws.Range("D12").PasteSpecial(-4122)

$ws.Range("A12").Value = "The QualOSS Open Source Assessment Model"
$ws.Range("B12").Value = "Martín Soto and Marcus Ciolkowski" + [char]10
$ws.Range("D12").Value = "https://www.rose-hulman.edu/class/csse/OldFiles/csse575/Resources/MeasOpenSource-05314237.pdf"
$ws.Range("C12").Value = 40087

$ws.Rows.Item(12).RowHeight = 45

# ---------------------------------------------------------------------------
# Row 13 (ht=45)
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A7").Copy()
$ws.Range("B13").PasteSpecial(-4122)

$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$ws.Range("A7").Copy()
$ws.Range("D13").PasteSpecial(-4122)

$ws.Range("A9").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E13").HorizontalAlignment = -4142
$ws.Range("E13").VerticalAlignment = -4160
$ws.Range("E13").WrapText = $false

$ws.Range("D13").Value = "http://dl.acm.org.libproxy.auc.ca/citation.cfm?id=1572200&CFID=951820277&CFTOKEN=35034244"
$ws.Range("A13").Value = "Introducing the OpenSource Maturity Model"
$ws.Range("B13").Value = "Etiel Petrinja, Ranga Nambakam, Alberto Sillitti"
$ws.Range("C13").Value = 39934

$ws.Rows.Item(13).RowHeight = 45

# ---------------------------------------------------------------------------
# Selection, as left by the editing session
# ---------------------------------------------------------------------------
$ws.Range("D23").Select()
